$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full set of rows (Sending cluster x Target cluster, 4x4) for the Gnai2-Tshr
# ligand-receptor pair table, reflecting the updated NATMI output values.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate, G Ligand average expression value,
# H Ligand total expression value, I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value, K Receptor-expressing cells,
# L Receptor detection rate, M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value, P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity, T Edge total expression derived specificity
$data = @(
    @("ECs", "Gnai2", "Tshr", "ECs", 3, 1, 150.0354306666667, 450.106292, 0.4152507364956075, 0.4152507364956075, 3, 1, 0.22859, 0.68577, 0.06730352972305123, 0.06730352972305123, 34.29659909609333, 308.66939186484, 0.02794784028625103, 0.02794784028625103),
    @("ECs", "Gnai2", "Tshr", "FAPs", 3, 1, 150.0354306666667, 450.106292, 0.4152507364956075, 0.4152507364956075, 3, 1, 2.180983666666667, 6.542951, 0.6421448840062525, 0.6421448840062525, 327.2248237052991, 2945.023413347692, 0.2666511360204828, 0.2666511360204827),
    @("ECs", "Gnai2", "Tshr", "M2", 3, 1, 150.0354306666667, 450.106292, 0.4152507364956075, 0.4152507364956075, 2, 0.6666666666666666, 0.173517, 0.520551, 0.05108844029465279, 0.05108844029465278, 26.033697822988, 234.303280406892, 0.02121451245876644, 0.02121451245876644),
    @("ECs", "Gnai2", "Tshr", "sCs", 3, 1, 150.0354306666667, 450.106292, 0.4152507364956075, 0.4152507364956075, 3, 1, 0.8133136666666667, 2.439941, 0.2394631459760435, 0.2394631459760435, 122.0258662454191, 1098.232796208772, 0.09943724773010724, 0.09943724773010722),
    @("FAPs", "Gnai2", "Tshr", "ECs", 3, 1, 68.382243, 205.146729, 0.1892604742946246, 0.1892604742946246, 3, 1, 0.22859, 0.68577, 0.06730352972305123, 0.06730352972305123, 15.63149692737, 140.68347234633, 0.01273789795708704, 0.01273789795708704),
    @("FAPs", "Gnai2", "Tshr", "FAPs", 3, 1, 68.382243, 205.146729, 0.1892604742946246, 0.1892604742946246, 3, 1, 2.180983666666667, 6.542951, 0.6421448840062525, 0.6421448840062525, 149.140555073031, 1342.264995657279, 0.1215326453128901, 0.12153264531289),
    @("FAPs", "Gnai2", "Tshr", "M2", 3, 1, 68.382243, 205.146729, 0.1892604742946246, 0.1892604742946246, 2, 0.6666666666666666, 0.173517, 0.520551, 0.05108844029465279, 0.05108844029465278, 11.865481658631, 106.789334927679, 0.009669022441138599, 0.009669022441138595),
    @("FAPs", "Gnai2", "Tshr", "sCs", 3, 1, 68.382243, 205.146729, 0.1892604742946246, 0.1892604742946246, 3, 1, 0.8133136666666667, 2.439941, 0.2394631459760435, 0.2394631459760435, 55.616212789221, 500.545915102989, 0.04532090858350893, 0.04532090858350891),
    @("M2", "Gnai2", "Tshr", "ECs", 3, 1, 104.737245, 314.211735, 0.2898796499701289, 0.2898796499701289, 3, 1, 0.22859, 0.68577, 0.06730352972305123, 0.06730352972305123, 23.94188683455, 215.47698151095, 0.01950992363787226, 0.01950992363787226),
    @("M2", "Gnai2", "Tshr", "FAPs", 3, 1, 104.737245, 314.211735, 0.2898796499701289, 0.2898796499701289, 3, 1, 2.180983666666667, 6.542951, 0.6421448840062525, 0.6421448840062525, 228.430220636665, 2055.871985729985, 0.1861447342058415, 0.1861447342058415),
    @("M2", "Gnai2", "Tshr", "M2", 3, 1, 104.737245, 314.211735, 0.2898796499701289, 0.2898796499701289, 2, 0.6666666666666666, 0.173517, 0.520551, 0.05108844029465279, 0.05108844029465278, 18.173692540665, 163.563232865985, 0.01480949919013378, 0.01480949919013378),
    @("M2", "Gnai2", "Tshr", "sCs", 3, 1, 104.737245, 314.211735, 0.2898796499701289, 0.2898796499701289, 3, 1, 0.8133136666666667, 2.439941, 0.2394631459760435, 0.2394631459760435, 85.184232767515, 766.658094907635, 0.06941549293628138, 0.06941549293628137),
    @("sCs", "Gnai2", "Tshr", "ECs", 3, 1, 38.15794, 114.47382, 0.105609139239639, 0.105609139239639, 3, 1, 0.22859, 0.68577, 0.06730352972305123, 0.06730352972305123, 8.7225235046, 78.5027115414, 0.0071078678418409, 0.007107867841840899),
    @("sCs", "Gnai2", "Tshr", "FAPs", 3, 1, 38.15794, 114.47382, 0.105609139239639, 0.105609139239639, 3, 1, 2.180983666666667, 6.542951, 0.6421448840062525, 0.6421448840062525, 83.22184389364668, 748.9965950428201, 0.06781636846703816, 0.06781636846703813),
    @("sCs", "Gnai2", "Tshr", "M2", 3, 1, 38.15794, 114.47382, 0.105609139239639, 0.105609139239639, 2, 0.6666666666666666, 0.173517, 0.520551, 0.05108844029465279, 0.05108844029465278, 6.621051274980001, 59.58946147482, 0.00539540620461397, 0.005395406204613968),
    @("sCs", "Gnai2", "Tshr", "sCs", 3, 1, 38.15794, 114.47382, 0.105609139239639, 0.105609139239639, 3, 1, 0.8133136666666667, 2.439941, 0.2394631459760435, 0.2394631459760435, 31.03437409384667, 279.30936684462, 0.02528949672614598, 0.02528949672614597)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
